# Auto-generated COM script applying the cryptos.xlsx price/volume update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.655.96'
$ws.Range('E2').Value = '  +4.21%  '
$ws.Range('D3').Value = '3.075.72'
$ws.Range('E3').Value = '  +5.54%  '
$ws.Range('E4').Value = '  -0.77%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '518.31'
$ws.Range('E5').Value = '  +4.55%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.33'
$ws.Range('E6').Value = '  +4.90%  '
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.434'
$ws.Range('E8').Value = '  +2.50%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '7.29'
$ws.Range('E9').Value = '  +4.43%  '
$ws.Range('E10').Value = '  +3.77%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.375'
$ws.Range('E11').Value = '  +4.09%  '
$ws.Range('D12').Value = '3.601.73'
$ws.Range('E12').Value = '  +2.42%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.129'
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.74'
$ws.Range('E14').Value = '  +0.18%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000164'
$ws.Range('E15').Value = '  +3.09%  '
$ws.Range('D16').Value = '57.685.04'
$ws.Range('E16').Value = '  +3.55%  '
$ws.Range('D17').Value = '3.075.54'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.09'
$ws.Range('E18').Value = '  +3.77%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.05'
$ws.Range('E19').Value = '  +2.42%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '8.12'
$ws.Range('E20').Value = '  +5.43%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '334.94'
$ws.Range('E21').Value = '  +4.56%  '
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.503'
$ws.Range('E23').Value = '  +3.75%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.99'
$ws.Range('E24').Value = '  +3.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.170'
$ws.Range('E25').Value = '  +7.04%  '
$ws.Range('E26').Value = '  -1.20%  '
$ws.Range('D27').Value = '0.0₃0917'
$ws.Range('E27').Value = '  +6.15%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.37'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.21'
$ws.Range('E29').Value = '  +5.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.83'
$ws.Range('E30').Value = '  +4.04%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.89'
$ws.Range('E31').Value = '  +5.30%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.18'
$ws.Range('E32').Value = '  +3.74%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '154.59'
$ws.Range('E33').Value = '  +3.04%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.48'
$ws.Range('E34').Value = '  -0.23%  '
$ws.Range('B35').Value = 'EnergySwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '27.08'
$ws.Range('E35').Value = '  +10.83%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.92'
$ws.Range('E36').Value = '  +4.26%  '
$ws.Range('E37').Value = '  +4.28%  '
$ws.Range('E38').Value = '  +4.45%  '
$ws.Range('D39').Value = '3.112.19'
$ws.Range('E39').Value = '  +3.25%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.93'
$ws.Range('E40').Value = '  +6.64%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.98'
$ws.Range('E41').Value = '  +1.50%  '
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.656'
$ws.Range('E43').Value = '  +2.18%  '
$ws.Range('D44').Value = '2.269.30'
$ws.Range('E44').Value = '  +5.79%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0257'
$ws.Range('E45').Value = '  +9.67%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.39'
$ws.Range('E46').Value = '  +4.75%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '20.43'
$ws.Range('E47').Value = '  +5.46%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.89'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.930'
$ws.Range('E49').Value = '  +1.88%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '265.15'
$ws.Range('E50').Value = '  +18.36%  '
$ws.Range('E51').Value = '  +4.22%  '
